# Applies the "Committing the latest code" change set:
#  - Adds 3 new worksheets (OrderShippingPageTest, OrderPaymentPageTest,
#    OrderConfirmationPageTest) after OrderDetailsPageTest, populated with
#    the new order-shipping/payment/confirmation test data.
#  - Updates selections on ShoppingCartPageTest / OrderDetailsPageTest.
#  - Leaves the newly-added OrderConfirmationPageTest tab active/selected.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) OrderShippingPageTest - duplicate of OrderDetailsPageTest's content
# ---------------------------------------------------------------------
$orderDetails = $wb.Worksheets.Item("OrderDetailsPageTest")
$orderDetails.Copy([System.Reflection.Missing]::Value, $orderDetails)
$shipping = $wb.Worksheets.Item($orderDetails.Index + 1)
$shipping.Name = "OrderShippingPageTest"
$shipping.Range("C14").Select()

# ---------------------------------------------------------------------
# 2) OrderPaymentPageTest - new sheet with payment-mode test data
# ---------------------------------------------------------------------
$payment = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $shipping)
$payment.Name = "OrderPaymentPageTest"

$payment.Range("A1").Value = "PaymentMode"
$payment.Range("B2").Value = "You have chosen to pay by check. Here is a short summary of your order:"
$payment.Range("B1").Value = "PaymentText"
$payment.Range("C1").Value = "TotalAmountOfProduct"
$payment.Range("A2").Value = "CHECK PAYMENT"
$payment.Range("C2").NumberFormat = "@"
$payment.Range("C2").Value = "`$18.51"

$payment.Range("A2:C2").NumberFormat = "@"

$payment.Columns.Item(1).ColumnWidth = 23.08984375
$payment.Columns.Item(2).ColumnWidth = 62.54296875
$payment.Columns.Item(3).ColumnWidth = 20.36328125

$payment.PageSetup.Orientation = 1
$payment.Cells.Select()

# ---------------------------------------------------------------------
# 3) OrderConfirmationPageTest - new sheet with confirmation test data
# ---------------------------------------------------------------------
$confirmation = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $payment)
$confirmation.Name = "OrderConfirmationPageTest"

$confirmation.Range("A2").Value = "ORDER CONFIRMATION"
$confirmation.Range("B2").Value = "Your order on My Store is complete."
$confirmation.Range("C2").Value = "Your order will be sent as soon as we receive your payment."
$confirmation.Range("A1").Value = "OrderConfimationPageLabel"
$confirmation.Range("B1").Value = "OrderConfimationSuccessMessage"
$confirmation.Range("C1").Value = "OrderConfirmationText"
$confirmation.Range("D1").Value = "TotalAmountOfProduct"
$confirmation.Range("D2").NumberFormat = "@"
$confirmation.Range("D2").Value = "`$18.51"

$confirmation.Range("A1:D2").NumberFormat = "@"

$confirmation.Columns.Item(1).ColumnWidth = 24.26953125
$confirmation.Columns.Item(2).ColumnWidth = 62.54296875
$confirmation.Columns.Item(3).ColumnWidth = 51.54296875
$confirmation.Columns.Item(4).ColumnWidth = 20.36328125

$confirmation.PageSetup.Orientation = 1
$confirmation.Range("B11").Select()

# ---------------------------------------------------------------------
# 4) Tweak selections on the pre-existing sheets touched by this commit
# ---------------------------------------------------------------------
$cart = $wb.Worksheets.Item("ShoppingCartPageTest")
$cart.Range("D2").Select()

$orderDetails.Cells.Select()

# ---------------------------------------------------------------------
# 5) Leave OrderConfirmationPageTest as the active/selected tab
# ---------------------------------------------------------------------
$confirmation.Activate()
$confirmation.Range("B11").Select()
